$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (correct marks per question) and "Total" row values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
